# NHAIRI_Planning_TPI.xlsx - add the "10.05.2017" tracking sheet
# (next day's copy of the planning sheet) with updated status markers.

$wb = $excel.ActiveWorkbook

# The existing sheet is the "09.05.2017" planning tab.
$day1 = $wb.Worksheets.Item(1)

# Duplicate it right after itself to create the next day's sheet.
$day1.Copy($null, $day1)
$day2 = $wb.Worksheets.Item(2)
$day2.Name = "10.05.2017"

# Row 8 ("Développement vue application C#" status markers):
#   - column C goes from blank to "Fait" (copy format+value from C5, an
#     existing "Fait" marker cell in the same style column).
#   - column D goes from "À faire" back to blank, so clear its contents
#     and restore the plain (unfilled) cell formatting from E8.
$day2.Range("C5").Copy($day2.Range("C8"))

$day2.Range("D8").ClearContents()
$day2.Range("E8").Copy()
$day2.Range("D8").PasteSpecial(-4122)

# Row 9 ("Développement modèles application C#"): column D becomes a new
# "À faire" marker, matching the existing one in E9.
$day2.Range("E9").Copy($day2.Range("D9"))

# Update the active selections: the old sheet keeps a plain selection,
# the new (now active) sheet is selected on E7.
$day1.Range("C8").Select() | Out-Null
$day2.Range("E7").Select() | Out-Null
$day2.Activate() | Out-Null
